# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:G values per row (column A - the date - and the header row are untouched).
# Column order per row: B, C, D, E, F, G
$values = @{
    2 = @(0.003994804209775715, 0.00007097389502863649, 0.8054896365839992, 0.496779210170732, 1, 1.306334624859536)
    3 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
    4 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    5 = @(0.04763786555579896, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 1.392311199053152)
    6 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    7 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    8 = @(0.6753301551942219, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 3.645393585217082)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($i = 0; $i -lt $rowValues.Count; $i++) {
        $col = 2 + $i   # column B = 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}
